$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.017.01'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -0.03%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.629.11'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.83%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.06'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.503'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.81%  '
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('E8').Value = '  -1.72%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0619'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -3.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.50'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -5.53%  '
$ws.Range('E11').Value = '  -0.90%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.856.77'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.78%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.634.05'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.14%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.18'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.96%  '
$ws.Range('E15').Value = '  -3.01%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '26.025.92'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.12%  '
$ws.Range('E17').Value = '  -2.35%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '61.48'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -3.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.01'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.06%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '192.62'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.08%  '
$ws.Range('E21').Value = '  -2.51%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.53'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -3.74%  '
$ws.Range('E23').Value = '  -2.26%  '
$ws.Range('E24').Value = '  +0.91%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.38'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.27%  '
$ws.Range('E26').Value = '  +0.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.73'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -3.27%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.71'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -2.70%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.24'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.72%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.23'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.75%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0482'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -2.60%  '
$ws.Range('E32').Value = '  -4.07%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.12'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -4.86%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.42'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.31%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.49'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -3.14%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.127.67'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.23%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.849'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -5.97%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.43'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.28%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.518'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -3.58%  '
$ws.Range('E40').Value = '  -2.17%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '98.16'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.66%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.766.36'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.83%  '
$ws.Range('E43').Value = '  -4.75%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.11'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -5.87%  '
$ws.Range('E45').Value = '  -1.27%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '54.39'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -3.72%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0524'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.44%  '
$ws.Range('E48').Value = '  -0.84%  '
$ws.Range('E49').Value = '  -0.36%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.48'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -3.83%  '
$ws.Range('E51').Value = '  +0.31%  '
